$d = $word.ActiveDocument

# 1. Update the letter date.
$d.Content.Find.Execute(
    "September 19, 2025", $true, $false, $false, $false, $false,
    $true, 1, $false, "September 21, 2025", 2) | Out-Null

# 2. Split the mailing address line "20635 Maria Court, Castro Valley CA 94546"
#    into two lines plus a trailing blank line:
#      20635 Maria Court
#      Castro Valley, CA 94546
#      <blank>
#    "^p" in the replacement text inserts a paragraph mark, so this both
#    re-words the text and creates the two new paragraphs in one step.
$d.Content.Find.Execute(
    "20635 Maria Court, Castro Valley CA 94546", $true, $false, $false,
    $false, $false, $false, 1, $false,
    "20635 Maria Court^pCastro Valley, CA 94546^p", 2) | Out-Null

# 3. Remove the two blank paragraphs (No Spacing, then Title) that sit
#    right after the "...Board of Directors" signature line, leaving the
#    remaining blank Title paragraph in place.
$found = $false
For ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.StartsWith("Lorena Circle Homeowners Association Board of Directors")) {
        $found = $true
        $d.Paragraphs.Item($i + 1).Range.Delete() | Out-Null
        $d.Paragraphs.Item($i + 1).Range.Delete() | Out-Null
        break
    }
}
